# Update "想去人数" (interested-people count) values in column F
# across the 展览 (sheet1), 演出 (sheet2) and 全部类型 (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1180
$ws1.Range("F3").Value = 640
$ws1.Range("F4").Value = 348
$ws1.Range("F6").Value = 523
$ws1.Range("F7").Value = 9298
$ws1.Range("F8").Value = 239
$ws1.Range("F9").Value = 526
$ws1.Range("F10").Value = 83
$ws1.Range("F11").Value = 645

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 20

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1180
$ws4.Range("F4").Value = 348
$ws4.Range("F5").Value = 14
$ws4.Range("F6").Value = 20
$ws4.Range("F7").Value = 4992
$ws4.Range("F8").Value = 0
$ws4.Range("F10").Value = 9298
$ws4.Range("F11").Value = 239
$ws4.Range("F12").Value = 526
$ws4.Range("F13").Value = 83
$ws4.Range("F14").Value = 6
$ws4.Range("F15").Value = 3
$ws4.Range("F16").Value = 645
$ws4.Range("F17").Value = 74
